$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 3.123541145015474
$ws.Range("E2").Value = 3.275761698145385
$ws.Range("C3").Value = 3.959010658874851
$ws.Range("E3").Value = 3.716811705074696
$ws.Range("C4").Value = 4.722695063536686
$ws.Range("E4").Value = 5.522497812297966
$ws.Range("C5").Value = 6.739021039846627
$ws.Range("E5").Value = 3.88772167754905
$ws.Range("C6").Value = 2.619839412265601
$ws.Range("E6").Value = 2.645362710332533
$ws.Range("C7").Value = -0.7919564768266385
$ws.Range("E7").Value = 0.8976572162586516
$ws.Range("C8").Value = 1.877689851450803
$ws.Range("E8").Value = 1.899791808163398
$ws.Range("C9").Value = 2.705004599189187
$ws.Range("E9").Value = 2.111643953433728
$ws.Range("C10").Value = 1.110374544249249
$ws.Range("E10").Value = 2.576021643263426
$ws.Range("C11").Value = 2.267566233338814
$ws.Range("E11").Value = 2.348556921565126
$ws.Range("C12").Value = 2.688433258834588
$ws.Range("E12").Value = 2.730428819177333
$ws.Range("C13").Value = 1.014079695989589
$ws.Range("E13").Value = 1.901826179618205
$ws.Range("C14").Value = 3.013853578092252
$ws.Range("E14").Value = 2.435357021275819
$ws.Range("C15").Value = 1.331295149770684
$ws.Range("E15").Value = 1.127665471558248
$ws.Range("C16").Value = 0.04589006555719699
$ws.Range("E16").Value = 0.8401596151991431
$ws.Range("C17").Value = 0.009546395482029624
$ws.Range("E17").Value = 0.02245646656315881
$ws.Range("C18").Value = 0.8709390141433015
$ws.Range("E18").Value = 0.7407001102931465
$ws.Range("C19").Value = 0.7652063367885598
$ws.Range("E19").Value = 1.397861196490657
$ws.Range("C20").Value = 2.267579219134386
$ws.Range("E20").Value = 2.220975586034668
$ws.Range("C21").Value = 3.146753122914103
$ws.Range("E21").Value = 2.456586080053058
$ws.Range("C22").Value = 1.769033835366818
$ws.Range("E22").Value = 1.086632508372576
$ws.Range("C23").Value = -4.774715709990263
$ws.Range("E23").Value = -0.8529145826070339
$ws.Range("C24").Value = 1.95493704440024
$ws.Range("E24").Value = 1.373265374526711
$ws.Range("C25").Value = 3.478075069442799
$ws.Range("E25").Value = 2.807281147895924
$ws.Range("C26").Value = 1.232342134690434
$ws.Range("E26").Value = 2.058767060726563
$ws.Range("C27").Value = 0.2542811494408159
$ws.Range("E27").Value = 1.100861823237564
$ws.Range("C28").Value = 1.519778766382096
$ws.Range("E28").Value = 0.6176362615319508
$ws.Range("C29").Value = 1.469441753880329
$ws.Range("E29").Value = 1.593307036690472
$ws.Range("C30").Value = 1.638203081492495
$ws.Range("E30").Value = 1.425381384581903
$ws.Range("C31").Value = 2.268697431234346
$ws.Range("E31").Value = 2.525738252590148
$ws.Range("C32").Value = 1.984425467899631
$ws.Range("E32").Value = 1.272991730262341
$ws.Range("C33").Value = 0.6066448776129052
$ws.Range("E33").Value = 0.8323605077169782
$ws.Range("C34").Value = -4.243076347305386
$ws.Range("E34").Value = -0.2172115314521883
$ws.Range("C35").Value = 1.438499295329754
$ws.Range("E35").Value = 1.725839624491665
$ws.Range("C36").Value = 1.906593537051537
$ws.Range("E36").Value = 1.222362529774923
$ws.Range("C37").Value = 0.08348019664223827
$ws.Range("E37").Value = 0.9116426337375527
$ws.Range("C38").Value = -0.214505326882275
$ws.Range("E38").Value = -0.02742999272021818
$ws.Range("C39").Value = 0.1651547428133782
$ws.Range("E39").Value = -0.1331392688890709
